$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.812.82'
$ws.Range("E2").Value = '  +1.72%  '
$ws.Range("D3").Value = '3.454.93'
$ws.Range("E3").Value = '  +0.59%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.86'
$ws.Range("E5").Value = '  -0.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.13'
$ws.Range("E6").Value = '  +2.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '3.455.34'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.583'
$ws.Range("E9").Value = '  +8.99%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.36'
$ws.Range("E10").Value = '  -2.83%  '
$ws.Range("E11").Value = '  +2.41%  '
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("D13").Value = '4.053.50'
$ws.Range("E13").Value = '  +0.76%  '
$ws.Range("E14").Value = '  -2.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000194'
$ws.Range("E15").Value = '  +4.04%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.12'
$ws.Range("E16").Value = '  +3.13%  '
$ws.Range("D17").Value = '64.883.18'
$ws.Range("E17").Value = '  +1.73%  '
$ws.Range("D18").Value = '3.460.84'
$ws.Range("E18").Value = '  +0.43%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.36'
$ws.Range("E19").Value = '  -1.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.29'
$ws.Range("E20").Value = '  -0.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '387.72'
$ws.Range("E21").Value = '  -1.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.19'
$ws.Range("E22").Value = '  -3.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.23'
$ws.Range("E23").Value = '  +1.58%  '
$ws.Range("E24").Value = '  +0.69%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("E26").Value = '  +15.77%  '
$ws.Range("E27").Value = '  +1.41%  '
$ws.Range("E28").Value = '  -0.72%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.24'
$ws.Range("E30").Value = '  +8.27%  '
$ws.Range("E31").Value = '  +4.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.04'
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.70'
$ws.Range("E33").Value = '  +0.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.53'
$ws.Range("E34").Value = '  -1.93%  '
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("E36").Value = '  +4.30%  '
$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '163.32'
$ws.Range("E37").Value = '  +2.71%  '
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.50'
$ws.Range("E38").Value = '  +0.42%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.91'
$ws.Range("E39").Value = '  +1.10%  '
$ws.Range("D40").Value = '3.015.80'
$ws.Range("E40").Value = '  +3.19%  '
$ws.Range("E41").Value = '  -2.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.26'
$ws.Range("E42").Value = '  -2.58%  '
$ws.Range("E43").Value = '  +4.85%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '42.94'
$ws.Range("E44").Value = '  +2.63%  '
$ws.Range("E45").Value = '  -1.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.775'
$ws.Range("E46").Value = '  +0.67%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.44'
$ws.Range("E47").Value = '  +7.98%  '
$ws.Range("E48").Value = '  +0.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.873'
$ws.Range("E49").Value = '  +5.66%  '
$ws.Range("E50").Value = '  +3.23%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.17'
$ws.Range("E51").Value = '  +6.04%  '